$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
